# Append the new BleepingComputer article as row 65 of the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A65").Value = "BleepingComputer"
$ws.Range("B65").Value = "Ongoing Duo outage causes Azure Auth authentication errors"
$ws.Range("C65").Value = "https://www.bleepingcomputer.com/news/technology/ongoing-duo-outage-causes-azure-auth-authentication-errors/"

# The publish date looks like a date literal ("2023-08-21"); Excel would
# normally auto-convert that into a date serial number on assignment. The
# source data stores it as plain text, so enter it with a leading
# apostrophe (text qualifier) and then strip the resulting "quote prefix"
# cell style so the cell ends up as plain, unstyled text - matching how
# every other row in this column is stored.
$ws.Range("D65").Value = "'2023-08-21"
$ws.Range("D65").Style = "Normal"

$ws.Range("E65").Value = "Cisco-owned multi-factor authentication (MFA) provider Duo Security is investigating an ongoing outage that has been causing authentication failures and errors starting three hours ago. [...]"
